$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Add a new row for the "Stock with txn fee DP" question, solved via
#    recursion + memoization. Insert it below row 3 FIRST (while row 3
#    still carries its original formatting), so the new row inherits
#    row 3's date / hyperlink-font styles.
# ------------------------------------------------------------------
$ws.Rows("4:4").Insert()              # inherits row 3's original formatting (date style, hyperlink font)

$ws.Range("A4").Value = [datetime]"2026-02-13"
$ws.Range("B4").Value = "Stock with txn fee DP"
$ws.Range("C4").Value = "https://leetcode.com/problems/best-time-to-buy-and-sell-stock-with-transaction-fee/"
$ws.Hyperlinks.Add($ws.Range("C4"), "https://leetcode.com/problems/best-time-to-buy-and-sell-stock-with-transaction-fee/")

# ------------------------------------------------------------------
# 2) Correct the date on the existing "Stock cooldown DP" row (row 3):
#    46093 (2026-03-12) -> 46065 (2026-02-12). Question/URL text stay
#    the same; only the date changes.
# ------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)   # xlPasteFormats - reuse the date style already used by A2
$ws.Range("A3").Value = [datetime]"2026-02-12"

# ------------------------------------------------------------------
# 3) Widen columns B and C to fit the new, longer text.
# ------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 18.15
$ws.Columns("C").ColumnWidth = 74

Write-Host "Applied tracking-question update for stock-with-transaction-fee."
